$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert the new rows (structural) ---
# Insert 1 row before old row 41 (new "CompletedFolder" row)
$ws.Rows("41:41").Insert()
# Insert 16 rows before what is now row 50 (old row 49, "Regex" section)
$ws.Rows("50:65").Insert()

# --- Step 2: Copy cell formatting from template rows/cells onto the new rows ---
# 3-column wrap-text data row template (row 35: CPHLinkLookUp) -> row 41 and row 50
$ws.Range("A35:C35").Copy()
$ws.Range("A41:C41").PasteSpecial(-4122)
$ws.Range("A35:C35").Copy()
$ws.Range("A50:C50").PasteSpecial(-4122)

# 2-column wrap-text data row template (row 43: SAGMailbox, A/B only) -> rows 51-60
$ws.Range("A43:B43").Copy()
$ws.Range("A51:B51").PasteSpecial(-4122)
$ws.Range("A52:B52").PasteSpecial(-4122)
$ws.Range("A53:B53").PasteSpecial(-4122)
$ws.Range("A54:B54").PasteSpecial(-4122)
$ws.Range("A55:B55").PasteSpecial(-4122)
$ws.Range("A56:B56").PasteSpecial(-4122)
$ws.Range("A57:B57").PasteSpecial(-4122)
$ws.Range("A58:B58").PasteSpecial(-4122)
$ws.Range("A59:B59").PasteSpecial(-4122)
$ws.Range("A60:B60").PasteSpecial(-4122)

# Single-column (A only) template (A43) -> A61:A65
$ws.Range("A43").Copy()
$ws.Range("A61:A65").PasteSpecial(-4122)
# Single wrap-text value cell template (B43) -> C61 (Name/Description-only row)
$ws.Range("B43").Copy()
$ws.Range("C61").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Step 3: Row heights for the two multi-line rows ---
$ws.Rows("41:41").RowHeight = 30
$ws.Rows("50:50").RowHeight = 180

# --- Step 4: Set cell values ---
# Row 41
$ws.Range("A41").Value = "CompletedFolder"
$ws.Range("B41").Value = "\\cavmfil001\Common\SinglePaymentScheme\Cross Compliance {0}\9)Processing Folders\Robot\{1}\Completed\"
$ws.Range("C41").Value = "Completed folder path for crfs processed successfully"

# Row 50
$ws.Range("A50").Value = "AttachmentNote"
$ws.Range("B50").Value = "Team Leader Name:  Rob Mclean`nInspection result: 2020 Cross Compliance {0} – {1}`n{2}`nVersion of the 2020 Cross Compliance Processing – {3}:`n{4}`nInspection result Letter/Email Sent: {5}`nCRF and Supporting Documents added to CRM: {6}"
$ws.Range("C50").Value = "Note to be added with last document attached in Customer Notification. `n"

# Row 51
$ws.Range("A51").Value = "SAGInspectionResult"
$ws.Range("B51").Value = "Sheep and Goat inspection"

# Row 52
$ws.Range("A52").Value = "RPAInspectionResult"
$ws.Range("B52").Value = "RPAi"

# Row 53
$ws.Range("A53").Value = "CIIInspectionResult"
$ws.Range("B53").Value = "Cattle Identification (CII) "

# Row 54
$ws.Range("A54").Value = "AWInspectionResult"
$ws.Range("B54").Value = "Animal Welfare inspection"

# Row 55
$ws.Range("A55").Value = "AHTBInspectionResult"
$ws.Range("B55").Value = "APHA TB Test"

# Row 56
$ws.Range("A56").Value = "SAGGuidanceTitle"
$ws.Range("B56").Value = "Sheep and Goat Inspection Instructions"

# Row 57
$ws.Range("A57").Value = "RPAGuidanceTitle"
$ws.Range("B57").Value = "Rural Payments Agency Inspection (RPAi) Instructions"

# Row 58
$ws.Range("A58").Value = "CIIGuidanceTitle"
$ws.Range("B58").Value = "CII Instructions used"

# Row 59
$ws.Range("A59").Value = "AWGuidanceTitle"
$ws.Range("B59").Value = "APHA Animal Welfare Inspection Instructions"

# Row 60
$ws.Range("A60").Value = "AHTBGuidanceTitle"
$ws.Range("B60").Value = "APHA Late TB Test Instructions"

# Row 61
$ws.Range("A61").Value = "SAGGuidanceVersion"
$ws.Range("C61").Value = "Version Number of Guidance"

# Row 62
$ws.Range("A62").Value = "RPAGuidanceVersion"

# Row 63
$ws.Range("A63").Value = "CIIGuidanceVersion"

# Row 64
$ws.Range("A64").Value = "AWGuidanceVersion"

# Row 65
$ws.Range("A65").Value = "AHTBGuidanceVersion"

# --- Step 5: Expand Table1 to the new range ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:C74"))

# --- Step 6: Update the sheet view (matches author scrolling to the new content) ---
$ws.Application.ActiveWindow.ScrollRow = 51
$ws.Range("B59").Select()
